$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Course Outcome"
$ws.Range("B1").Value = "Assignment1"
$ws.Range("C1").Value = "Assignment2"
$ws.Range("A2").Select()
